$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update forecast coefficient values in column B
$ws.Range("B2").Value = -0.5992774864336459
$ws.Range("B3").Value = 0.842982803180619
$ws.Range("B4").Value = 32.460182750837916

# Remove the last row (row 5), which also removes the now-unused "4" shared string
$ws.Rows("5:5").Delete()
